# Weekly price update: insert a new week's price record for "Apio"
# (Terminal Hortofrutícola Agro Chillán) ahead of the existing row 152,
# pushing the old row 152 (and everything after it) down by one row.
#
# The new row duplicates all the fixed/descriptive columns of the former
# row 152 (market, region, product, unit, etc.) and only differs in the
# date (column D) and the "cantidad" / volume figure (column J), which
# carry this week's new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 152 and insert the copy above it, shifting rows
# 152..178 down to 153..179 (dimension grows from R178 to R179).
$ws.Rows.Item(152).Copy()
$ws.Rows.Item(152).Insert()
$excel.CutCopyMode = $false

# Overwrite the new row's date and volume with this week's values.
$ws.Range("D152").Value = 44617
$ws.Range("J152").Value = 60
